# Actualización desde MV -datos-
# Adds the new daily volatility-index rows (MOVE / VIX) for
# 27-09-2021 .. 04-10-2021, matching the source workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 187 (27-09-2021) was already present but only had the VIX (C) value.
# Fill in the missing MOVE (B) value and correct the VIX (C) value.
$ws.Range("B187").Value = 60.5
$ws.Range("C187").Value = 18.76

# New row: 28-09-2021
$ws.Range("A188").Value = "28-09-2021"
$ws.Range("B188").Value = 62.71
$ws.Range("C188").Value = 23.25

# New row: 29-09-2021
$ws.Range("A189").Value = "29-09-2021"
$ws.Range("B189").Value = 61.21
$ws.Range("C189").Value = 22.56

# New row: 30-09-2021
$ws.Range("A190").Value = "30-09-2021"
$ws.Range("B190").Value = 61.07
$ws.Range("C190").Value = 23.14

# New row: 01-10-2021
# Leading apostrophe forces the date-looking label to stay plain text
# (otherwise it gets auto-converted into a date serial number); the
# Style reset afterwards keeps the cell on the workbook's default style
# instead of picking up a quote-prefix number format.
$ws.Range("A191").Value = "'01-10-2021"
$ws.Range("A191").Style = "Normal"
$ws.Range("B191").Value = 57.38
$ws.Range("C191").Value = 21.15

# New row: 04-10-2021 (no MOVE value published yet, only VIX)
$ws.Range("A192").Value = "'04-10-2021"
$ws.Range("A192").Style = "Normal"
$ws.Range("C192").Value = 22.98
